$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting the existing "Nombre" (and
# everything after it) one column to the right, to make room for the new
# "Clave" field.
$ws.Columns("C:C").Insert()

# Copy the header formatting from the neighboring "Usuario" header (B4) onto
# the new header cell so it matches the rest of the header row, then set its
# text to "Clave".
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C4").Value = "Clave"

# Match the resulting selection/active cell recorded in the saved view.
$ws.Range("C4").Select()
